$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the stray "_GoBack" bookmark that sits between the runs
#    "5" and " min" in the "15 min" table cell (near the top of the
#    document).  Bookmark elements are not part of the run text, so a
#    Range.Text assignment that is a no-op (same text in, same text
#    out) leaves them untouched; we briefly replace the 2 characters
#    spanning the bookmark with a placeholder and then restore the
#    original text so Word rebuilds that bit of the paragraph without
#    the bookmark.
# ------------------------------------------------------------------
$timeRng = $d.Content
$timeRng.Find.Execute("15 min", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $timeRng.Find.Found) {
    throw "Could not find '15 min' cell"
}
$bmSpan = $d.Range($timeRng.Start + 1, $timeRng.Start + 3)
$bmSpan.Text = "#~"
$bmSpan = $d.Range($timeRng.Start + 1, $timeRng.Start + 3)
$bmSpan.Text = "5 "

# ------------------------------------------------------------------
# 2) Reword the red italic bullet "4-state models?" to "3-state
#    model".
# ------------------------------------------------------------------
$stateRng = $d.Content
$stateRng.Find.Execute("4-state models?", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $stateRng.Find.Found) {
    throw "Could not find '4-state models?' bullet"
}
$stateRng.Text = "3-"
$tail = $d.Range($stateRng.End, $stateRng.End)
$tail.InsertAfter("state model")

# ------------------------------------------------------------------
# 3) Re-insert the "_GoBack" bookmark so it now wraps "modify models
#    based on audience input" in the following bullet.
# ------------------------------------------------------------------
$targetRng = $d.Content
$targetRng.Find.Execute("modify models based on audience input", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $targetRng.Find.Found) {
    throw "Could not find 'modify models based on audience input' text"
}
$d.Bookmarks.Add("_GoBack", $targetRng)
